# Applies the "solve errors of the example excels" edit:
#  - Sheet1 "table attribute": drop the 2-column attribute-name/value table and
#    replace it with a single column containing "table name " / "packet id allocate".
#  - Sheet2 "table content col attribute": rename header B1 to "column type", add a
#    new header C1 "column description" (with a new column width), becomes the
#    active/selected sheet.
#  - Sheet3 "table content": values stay the same (string table just gets reindexed
#    automatically by the engine).
#  - Workbook: active tab becomes sheet2 (index 1, 0-based).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Header (row 1) font color used throughout this workbook - sample it before
# we destroy any cells, so the new layout keeps the same look.
$headerColor = $ws2.Range("A1").Font.Color

# --- Sheet1 ("table attribute") ---
# Wipe the whole 2-column attribute-name/value table ...
$ws1.UsedRange.Delete()

# ... and rebuild it as a single column: table name / packet id allocate
$ws1.Range("A1").Value = "table name "
$ws1.Range("A1").Font.Color = $headerColor
$ws1.Range("A2").Value = "packet id allocate"

$ws1.Range("A2").Select() | Out-Null

# --- Sheet2 ("table content col attribute") ---
$ws2.Range("B1").Value = "column type"
$ws2.Range("C1").Value = "column description"

# Match the header look (green font) used by the rest of row 1
$ws2.Range("C1").Font.Color = $headerColor

# New column width for the added column C
$ws2.Columns.Item(3).ColumnWidth = 23

$ws2.Range("C4").Select() | Out-Null

# --- Sheet3 ("table content") : values unchanged, nothing to do ---

# --- Workbook: make sheet2 the active / selected tab ---
$ws2.Activate()
